# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310   (source/"before" format version)
#   *_new -> *_FV2404   (target/"after" format version)
# Then (re)create the table over the data range and freeze the header row,
# matching the regenerated AHB-diff export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<name>_old" -> "<name>_FV2310",
#                            "<name>_new" -> "<name>_FV2404"
# ---------------------------------------------------------------------------
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column  # xlToLeft
if ($lastCol -lt 1) { $lastCol = 21 }

for ($c = 1; $c -le $lastCol; $c++) {
    $headerCell = $ws.Cells.Item(1, $c)
    $headerText = $headerCell.Value2
    if ($null -eq $headerText) { continue }
    if (-not ($headerText -is [string])) { continue }

    if ($headerText.EndsWith("_old")) {
        $headerCell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2310"
    } elseif ($headerText.EndsWith("_new")) {
        $headerCell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2404"
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an actual Excel Table ("ListObject") so the
#    header row gets the table's autofilter + styling, matching Table1.
# ---------------------------------------------------------------------------
$tableRange = $ws.UsedRange

$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A1").Select()
